$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new note row (C19), italic font + yellow fill, matching the
# style used elsewhere in the sheet for the "class timing" note.
$cell = $ws.Range("C19")
$cell.Value = "* Class to review the work every Sunday @10AM "
$cell.Font.Italic = $true
$cell.Interior.Color = 65535

$ws.PageSetup.Orientation = 1

$ws.Range("C22").Select()
